# Add weekly Espárragos ("asparagus") price-report rows 22-27 to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mercado  = "Mercado Mayorista Lo Valledor de Santiago"
$region   = "Metropolitana"
$fecha    = 44476
$codreg   = 13
$catId    = 300000000
$categoria = "Espárragos"
$variedad  = "Sin especificar"
$unidad    = "`$/kilo"
$clasif    = "Hortaliza"

$rows = @(
    # row, calidad,   volumen, precioMin, precioMax, precioProm, origen,                     precioKg, kgUnid
    @(22, "Banquete", 1300,    1500,      1600,      1554,       "Provincia de Linares",      1554,     1),
    @(23, "Banquete", 700,     1400,      1500,      1457,       "Región Metropolitana",      1457,     1),
    @(24, "Primera",  900,     1300,      1400,      1356,       "Provincia de Linares",      1356,     1),
    @(25, "Primera",  500,     1200,      1300,      1260,       "Región Metropolitana",      1260,     1),
    @(26, "Segunda",  500,     1100,      1200,      1160,       "Provincia de Linares",      1160,     1),
    @(27, "Segunda",  200,     1000,      1100,      1050,       "Región Metropolitana",      1050,     1)
)

foreach ($r in $rows) {
    $rowIndex  = $r[0]
    $calidad   = $r[1]
    $volumen   = $r[2]
    $precioMin = $r[3]
    $precioMax = $r[4]
    $precioProm = $r[5]
    $origen    = $r[6]
    $precioKg  = $r[7]
    $kgUnid    = $r[8]

    $ws.Cells.Item($rowIndex, 1).Value  = 6
    $ws.Cells.Item($rowIndex, 2).Value  = $mercado
    $ws.Cells.Item($rowIndex, 3).Value  = $region
    $ws.Cells.Item($rowIndex, 4).Value  = $fecha
    $ws.Cells.Item($rowIndex, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($rowIndex, 5).Value  = $codreg
    $ws.Cells.Item($rowIndex, 6).Value  = $catId
    $ws.Cells.Item($rowIndex, 7).Value  = $categoria
    $ws.Cells.Item($rowIndex, 8).Value  = $variedad
    $ws.Cells.Item($rowIndex, 9).Value  = $calidad
    $ws.Cells.Item($rowIndex, 10).Value = $volumen
    $ws.Cells.Item($rowIndex, 11).Value = $precioMin
    $ws.Cells.Item($rowIndex, 12).Value = $precioMax
    $ws.Cells.Item($rowIndex, 13).Value = $precioProm
    $ws.Cells.Item($rowIndex, 14).Value = $unidad
    $ws.Cells.Item($rowIndex, 15).Value = $origen
    $ws.Cells.Item($rowIndex, 16).Value = $precioKg
    $ws.Cells.Item($rowIndex, 17).Value = $kgUnid
    $ws.Cells.Item($rowIndex, 18).Value = $clasif
}

Write-Output "Added rows 22-27"
